# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Rebuilds the "Periodo Mora" table (rows 16-46) on Hoja1: the old
# chronological ordering (oldest -> newest, grouped by worker) is replaced
# by a newest -> oldest ordering, and a batch of new periods (1801-1810,
# 1902) is appended for each worker while the G (Salario Basico) value for
# LUIS ALFREDO HERRERA ARRIETA is updated to 781242.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worker 1: MAURICIO VEGA ANAYA (CC 9203477) -> rows 16-36 ------------
$mauricio = @(
    @("1902", 28124, 689455),
    @("1810", 31249, 689455),
    @("1809", 31249, 689455),
    @("1808", 27578, 689455),
    @("1806", 27578, 689455),
    @("1804", 27578, 689455),
    @("1803", 27578, 689455),
    @("1802", 27578, 689455),
    @("1801", 27578, 689455),
    @("1712", 27578, 689455),
    @("1711", 27578, 689455),
    @("1710", 27578, 689455),
    @("1709", 27578, 689455),
    @("1708", 27578, 689455),
    @("1707", 27578, 689455),
    @("1706", 27578, 689455),
    @("1705", 27578, 689455),
    @("1704", 27578, 689455),
    @("1703", 27578, 689455),
    @("1702", 27578, 689455),
    @("1701", 27578, 689455)
)

# --- Worker 2: LUIS ALFREDO HERRERA ARRIETA (CC 73577314) -> rows 37-46 --
$luis = @(
    @("1902", 28124, 781242),
    @("1810", 31249, 781242),
    @("1809", 31249, 781242),
    @("1808", 29509, 781242),
    @("1806", 29509, 781242),
    @("1804", 29509, 781242),
    @("1803", 29509, 781242),
    @("1802", 29509, 781242),
    @("1801", 29509, 781242),
    @("1712", 29509, 781242)
)

$row = 16
foreach ($entry in $mauricio) {
    $ws.Range("B" + $row).Value = "CC"
    $ws.Range("C" + $row).Value = "9203477"
    $ws.Range("D" + $row).Value = "MAURICIO VEGA ANAYA"
    $ws.Range("E" + $row).Value = $entry[0]
    $ws.Range("F" + $row).Value = $entry[1]
    $ws.Range("G" + $row).Value = $entry[2]
    $row = $row + 1
}

foreach ($entry in $luis) {
    $ws.Range("B" + $row).Value = "CC"
    $ws.Range("C" + $row).Value = "73577314"
    $ws.Range("D" + $row).Value = "LUIS ALFREDO HERRERA ARRIETA"
    $ws.Range("E" + $row).Value = $entry[0]
    $ws.Range("F" + $row).Value = $entry[1]
    $ws.Range("G" + $row).Value = $entry[2]
    $row = $row + 1
}

# --- Column widths (best-fit columns were recalculated by Excel) --------
$ws.Columns.Item(2).ColumnWidth = 17.7109375
$ws.Columns.Item(3).ColumnWidth = 15.8935546875
$ws.Columns.Item(5).ColumnWidth = 12.7109375
$ws.Columns.Item(6).ColumnWidth = 9.3466796875
$ws.Columns.Item(7).ColumnWidth = 13.5302734375
$ws.Columns.Item(8).ColumnWidth = 18.5302734375
$ws.Columns.Item(9).ColumnWidth = 17.2568359375
$ws.Columns.Item(10).ColumnWidth = 14.1669921875
